# "make use of state table to simplify special handlers"
#
# This adds a new WINDUP -> SHOOT -> SOLENOID -> MENU sequence to the
# RubberBandGun state table, moving the "mSPCL_HANDLER | mSPCL_HANDLER_SHOOT"
# special-handler marker from row 4 down to its own new row (14), and adding
# a brand-new row (16) for the solenoid release handler, plus a little
# trailing formatted-but-empty row (17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 4: the "mSPCL_HANDLER | mSPCL_HANDLER_SHOOT" marker that used to
# live in B4 now moves down to row 14 (see below), so clear it here.
# I4 used to just repeat "mROW_MENU"; now it points at the new
# "mROW_SHOOT_WINDUP" row instead, and loses its old centering style.
# ---------------------------------------------------------------------
$ws.Range("B4").Clear()
$ws.Range("I4").Clear()
$ws.Range("I4").Value = "mROW_SHOOT_WINDUP"

# ---------------------------------------------------------------------
# Row 12: new "...and the WINDUP" state - winds up the shooter, then
# goes on to the new SHOOT row.
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "… and the WINDUP"
$ws.Range("C12").Value = "mROW_SHOOT_WINDUP"
$ws.Range("D12").Value = "mEFCT_WIND_UP"
$ws.Range("E12").Value = "mEFCT_WIND_UP"
$ws.Range("J12").Value = "mROW_SHOOT"

# ---------------------------------------------------------------------
# Row 14: new "POW!!!" SHOOT state, carrying the special handler marker
# that used to be on B4, then continuing on to the new SOLENOID row.
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "POW!!!"
$ws.Range("B14").VerticalAlignment = -4108
$ws.Range("B14").WrapText = $true
$ws.Range("B14").Value = "mSPCL_HANDLER | mSPCL_HANDLER_SHOOT"
$ws.Range("C14").Value = "mROW_SHOOT"
$ws.Range("D14").Value = "mEFCT_SHOOT"
$ws.Range("E14").Value = "mEFCT_SHOOT"
$ws.Range("J14").Value = "mROW_SOLENOID"

# ---------------------------------------------------------------------
# Row 16: new "release solenoid" state with its own special handler
# marker, silences the effect and returns to the MENU row.
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "release solenoid"
$ws.Range("B16").VerticalAlignment = -4108
$ws.Range("B16").WrapText = $true
$ws.Range("B16").Value = "mSPCL_HANDLER | mSPCL_HANDLER_SOLENOID"
$ws.Range("C16").Value = "mROW_SOLENOID"
$ws.Range("D16").Value = "mEFCT_UNIQ_SILENCE"
$ws.Range("E16").Value = "mEFCT_UNIQ_SILENCE"
$ws.Range("J16").VerticalAlignment = -4108
$ws.Range("J16").Value = "mROW_MENU"

# Row 14 and 16 both wrap long text in column B, making Excel grow those
# rows to a 30pt height (same as rows 4-6 above).
$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30

# ---------------------------------------------------------------------
# Row 17: trailing row, formatted like the others but left blank.
# ---------------------------------------------------------------------
$ws.Range("B17").VerticalAlignment = -4108
$ws.Range("B17").WrapText = $true

# ---------------------------------------------------------------------
# Column B grew slightly wider to fit the new, slightly longer text.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 25.75

# ---------------------------------------------------------------------
# Selection cosmetically moved to the newly-added E16 cell.
# ---------------------------------------------------------------------
$ws.Range("E16").Select() | Out-Null
